$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated cryptocurrency market data (prices, volume %, and two swapped-row identities)

# Row 2
$ws.Range('D2').Value = '60.479.89'
$ws.Range('E2').Value = '  +2.57%  '

# Row 3
$ws.Range('D3').Value = '2.705.41'
$ws.Range('E3').Value = '  +3.13%  '

# Row 4
$ws.Range('E4').Value = '  +0.14%  '

# Row 5
$c = $ws.Range('D5')
$c.Value = "'526.27"
$c.Style = 'Normal'
$ws.Range('E5').Value = '  +1.52%  '

# Row 6
$c = $ws.Range('D6')
$c.Value = "'145.00"
$c.Style = 'Normal'
$ws.Range('E6').Value = '  -0.04%  '

# Row 7
$c = $ws.Range('D7')
$c.Value = "'0.997"
$c.Style = 'Normal'
$ws.Range('E7').Value = '  +0.08%  '

# Row 8
$c = $ws.Range('D8')
$c.Value = "'0.576"
$c.Style = 'Normal'
$ws.Range('E8').Value = '  +1.99%  '

# Row 9
$ws.Range('D9').Value = '2.731.33'
$ws.Range('E9').Value = '  +3.18%  '

# Row 10
$c = $ws.Range('D10')
$c.Value = "'6.65"
$c.Style = 'Normal'
$ws.Range('E10').Value = '  +5.63%  '

# Row 11
$ws.Range('E11').Value = '  +1.22%  '

# Row 12
$ws.Range('E12').Value = '  +0.94%  '

# Row 13
$ws.Range('E13').Value = '  +3.03%  '

# Row 14
$ws.Range('D14').Value = '3.181.29'
$ws.Range('E14').Value = '  +2.93%  '

# Row 15
$ws.Range('D15').Value = '60.505.66'
$ws.Range('E15').Value = '  +2.69%  '

# Row 16
$ws.Range('B16').Value = 'Avalanche'
$ws.Range('C16').Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$c = $ws.Range('D16')
$c.Value = "'21.30"
$c.Style = 'Normal'
$ws.Range('E16').Value = '  +1.60%  '

# Row 17
$ws.Range('B17').Value = 'WrappedEther'
$ws.Range('C17').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D17').Value = '2.725.83'
$ws.Range('E17').Value = '  +3.04%  '

# Row 18
$ws.Range('E18').Value = '  +0.61%  '

# Row 19
$c = $ws.Range('D19')
$c.Value = "'347.90"
$c.Style = 'Normal'
$ws.Range('E19').Value = '  -0.63%  '

# Row 20
$ws.Range('E20').Value = '  -0.04%  '

# Row 21
$c = $ws.Range('D21')
$c.Value = "'10.62"
$c.Style = 'Normal'
$ws.Range('E21').Value = '  +3.08%  '

# Row 22
$c = $ws.Range('D22')
$c.Value = "'6.45"
$c.Style = 'Normal'
$ws.Range('E22').Value = '  +4.71%  '

# Row 23
$ws.Range('E23').Value = '  +0.04%  '

# Row 24
$ws.Range('E24').Value = '  +3.35%  '

# Row 25
$ws.Range('E25').Value = '  +0.55%  '

# Row 26
$ws.Range('E26').Value = '  +4.91%  '

# Row 27
$c = $ws.Range('D27')
$c.Value = "'0.994"
$c.Style = 'Normal'
$ws.Range('E27').Value = '  -0.09%  '

# Row 28
$ws.Range('D28').Value = '0.0₃0819'
$ws.Range('E28').Value = '  +2.05%  '

# Row 29
$ws.Range('E29').Value = '  +2.79%  '

# Row 30
$c = $ws.Range('D30')
$c.Value = "'6.81"
$c.Style = 'Normal'
$ws.Range('E30').Value = '  +9.24%  '

# Row 31
$ws.Range('E31').Value = '  +0.02%  '

# Row 32
$ws.Range('E32').Value = '  +1.48%  '

# Row 33
$c = $ws.Range('D33')
$c.Value = "'19.10"
$c.Style = 'Normal'
$ws.Range('E33').Value = '  +0.76%  '

# Row 34
$c = $ws.Range('D34')
$c.Value = "'150.48"
$c.Style = 'Normal'
$ws.Range('E34').Value = '  +0.51%  '

# Row 35
$c = $ws.Range('D35')
$c.Value = "'4.25"
$c.Style = 'Normal'
$ws.Range('E35').Value = '  +6.22%  '

# Row 36
$ws.Range('B36').Value = 'ImmutableX'
$ws.Range('C36').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$c = $ws.Range('D36')
$c.Value = "'1.23"
$c.Style = 'Normal'
$ws.Range('E36').Value = '  +8.21%  '

# Row 37
$ws.Range('B37').Value = 'SuiNetwork'
$ws.Range('C37').Value = 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
$c = $ws.Range('D37')
$c.Value = "'0.940"
$c.Style = 'Normal'
$ws.Range('E37').Value = '  -2.47%  '

# Row 38
$c = $ws.Range('D38')
$c.Value = "'0.874"
$c.Style = 'Normal'
$ws.Range('E38').Value = '  +4.02%  '

# Row 39
$ws.Range('E39').Value = '  +7.48%  '

# Row 40
$ws.Range('E40').Value = '  +0.97%  '

# Row 41
$ws.Range('E41').Value = '  -0.35%  '

# Row 42
$c = $ws.Range('D42')
$c.Value = "'283.32"
$c.Style = 'Normal'
$ws.Range('E42').Value = '  +2.50%  '

# Row 43
$c = $ws.Range('D43')
$c.Value = "'20.03"
$c.Style = 'Normal'
$ws.Range('E43').Value = '  +2.27%  '

# Row 44
$ws.Range('B44').Value = 'Mantle'
$ws.Range('C44').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$c = $ws.Range('D44')
$c.Value = "'0.612"
$c.Style = 'Normal'
$ws.Range('E44').Value = '  +0.61%  '

# Row 45
$ws.Range('B45').Value = 'Stellar'
$ws.Range('C45').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$c = $ws.Range('D45')
$c.Value = "'0.0988"
$c.Style = 'Normal'
$ws.Range('E45').Value = '  +0.42%  '

# Row 46
$ws.Range('B46').Value = 'FirstDigitalUSD'
$ws.Range('C46').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$c = $ws.Range('D46')
$c.Value = "'0.996"
$c.Style = 'Normal'
$ws.Range('E46').Value = '  +0.09%  '

# Row 47
$ws.Range('D47').Value = '2.143.23'
$ws.Range('E47').Value = '  +8.12%  '

# Row 48
$ws.Range('E48').Value = '  +3.16%  '

# Row 49
$ws.Range('B49').Value = 'RenderToken'
$ws.Range('C49').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$c = $ws.Range('D49')
$c.Value = "'4.82"
$c.Style = 'Normal'
$ws.Range('E49').Value = '  +2.18%  '

# Row 50
$c = $ws.Range('D50')
$c.Value = "'10.48"
$c.Style = 'Normal'
$ws.Range('E50').Value = '  +1.96%  '

# Row 51
$ws.Range('B51').Value = 'VeChain'
$ws.Range('C51').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$c = $ws.Range('D51')
$c.Value = "'0.0234"
$c.Style = 'Normal'
$ws.Range('E51').Value = '  +2.05%  '
